# Update background dataset: populate the GWP-LULUC_mfg (column K) values
# for all data rows with 0, matching the newly added column of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill K2:K288 with 0 (new GWP-LULUC_mfg values for every material row)
$ws.Range("K2:K288").Value = 0

# Reflect the final selection/view state used while editing the data
$ws.Range("K2:K288").Select()
